$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

# Update patient number and pincode values for the PlaceMedicineOrder flow
$ws1.Range("A2").Value = 6379948639
$ws1.Range("B2").Value = 600003

# Move to the payments tab area (selection moves to I1)
$ws1.Range("I1").Select()
